# "Minor changes to selectional Present absent panel."
# 1. Update a few values on Sheet1.
# 2. Add a new Sheet2 with its own small table.
# 3. Make Sheet2 the active sheet (activeTab goes from 0 -> 1).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("B2").Value = 16.0
$ws1.Range("B4").Value = 12.0
$ws1.Range("B5").Value = 0

# Insert the new sheet right after Sheet1 so final order is Sheet1, Sheet2.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$ws2.Range("A1").Value = "name"
$ws2.Range("B1").Value = "jan"
$ws2.Range("C1").Value = "feb"
$ws2.Range("D1").Value = "yello"

$ws2.Range("A2").Value = "kanjus"
$ws2.Range("B2").Value = 1.0

$ws2.Range("A3").Value = "maru"
$ws2.Range("B3").Value = 1.0

$ws2.Range("A4").Value = "laila"
$ws2.Range("B4").Value = 0

# Sheet2 becomes the active/selected tab.
$ws2.Activate()
